# Update the "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to match the values published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F
$updates = @{
    2  = 274
    3  = 1392
    4  = 163
    6  = 237
    8  = 17
    9  = 189
    11 = 4674
    12 = 6940
    16 = 576
    18 = 4161
    19 = 841
    21 = 66
    22 = 2741
    24 = 551
    25 = 174
    26 = 379
    27 = 377
    29 = 238
    30 = 45
    31 = 1644
    32 = 1041
    34 = 460
    36 = 551
    37 = 2
    38 = 497
    41 = 194
    42 = 648
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
